$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'basketball leg sleeve youth boys'
$ws.Cells.Item(2, 1).Value = 'youth basketball leg sleeves boys'
$ws.Cells.Item(3, 1).Value = 'basketball pants with knee pads'
$ws.Cells.Item(4, 1).Value = 'mens compression pants'
$ws.Cells.Item(5, 1).Value = 'compression knee pad'
$ws.Cells.Item(6, 1).Value = 'basketball leggings'
$ws.Cells.Item(7, 1).Value = 'knee compression pads'
$ws.Cells.Item(8, 1).Value = 'compression pants men'
$ws.Cells.Item(9, 1).Value = 'knee pads hex'
$ws.Cells.Item(10, 1).Value = 'basketball knee pads youth boys'
$ws.Cells.Item(11, 1).Value = 'mens capri pants'
$ws.Cells.Item(12, 1).Value = 'volleyball knee'
$ws.Cells.Item(13, 1).Value = 'workout hand pads'
$ws.Cells.Item(14, 1).Value = 'black athletic leggings'
$ws.Cells.Item(15, 1).Value = 'man pads'
$ws.Cells.Item(16, 1).Value = 'compression tights'
$ws.Cells.Item(17, 1).Value = 'tactical pants with knee pads'
$ws.Cells.Item(18, 1).Value = 'mens pants'
$ws.Cells.Item(19, 1).Value = 'xxl knee pads'
$ws.Cells.Item(20, 1).Value = 'knee compression for men'
$ws.Cells.Item(21, 1).Value = 'athletic leggings capri'
$ws.Cells.Item(22, 1).Value = 'sliding pads for softball'
$ws.Cells.Item(23, 1).Value = 'knee pads adult'
$ws.Cells.Item(24, 1).Value = 'softball knee pads'
$ws.Cells.Item(25, 1).Value = 'basketball knee pads adult'
$ws.Cells.Item(26, 1).Value = 'youth basketball gear'
$ws.Cells.Item(27, 1).Value = 'boys leggings'
$ws.Cells.Item(28, 1).Value = 'compression knee pads'
$ws.Cells.Item(29, 1).Value = 'black basketball leggings for men'
$ws.Cells.Item(30, 1).Value = 'compression pants with knee pads boys'
$ws.Cells.Item(31, 1).Value = 'hex knee pads basketball youth'
$ws.Cells.Item(32, 1).Value = 'anti strip clothing for men'
$ws.Cells.Item(33, 1).Value = 'hex protective knee pads'
$ws.Cells.Item(34, 1).Value = 'men leggings'
$ws.Cells.Item(35, 1).Value = 'weightlifting pads'
$ws.Cells.Item(36, 1).Value = 'compression tights boys'
$ws.Cells.Item(37, 1).Value = 'compression tights for boys'
$ws.Cells.Item(38, 1).Value = 'bjj tights for men'
$ws.Cells.Item(39, 1).Value = 'mens football pants'
$ws.Cells.Item(40, 1).Value = 'advanced squat pad'
$ws.Cells.Item(41, 1).Value = 'knee pads football'
$ws.Cells.Item(42, 1).Value = 'knee protector pain'
$ws.Cells.Item(43, 1).Value = 'knee support leggings'
$ws.Cells.Item(44, 1).Value = 'compression tights for youth'
$ws.Cells.Item(45, 1).Value = 'mens tall pants'
$ws.Cells.Item(46, 1).Value = 'men tights'
$ws.Cells.Item(47, 1).Value = 'wrestling gear'
$ws.Cells.Item(48, 1).Value = 'little boys athletic leggings'
$ws.Cells.Item(49, 1).Value = 'workout leggings for men pack'
$ws.Cells.Item(50, 1).Value = 'baseball gear for boys'
$ws.Cells.Item(51, 1).Value = 'black baseball pants youth xl'
$ws.Cells.Item(52, 1).Value = 'wrestling knee pad'
$ws.Cells.Item(53, 1).Value = 'volleyball knee pads xxl mens'
$ws.Cells.Item(54, 1).Value = 'knee running'
$ws.Cells.Item(55, 1).Value = 'knee pads for gym'
$ws.Cells.Item(56, 1).Value = 'knee for running'
$ws.Cells.Item(57, 1).Value = 'athletic compression leggings'
$ws.Cells.Item(58, 1).Value = 'weightlifting pad'
$ws.Cells.Item(59, 1).Value = 'foam knee pad'
$ws.Cells.Item(60, 1).Value = 'compression tight'
$ws.Cells.Item(61, 1).Value = 'knee support volleyball'
$ws.Cells.Item(62, 1).Value = 'basketball equipment'
$ws.Cells.Item(63, 1).Value = 'fitness squat pad'
$ws.Cells.Item(64, 1).Value = 'knee swelling'
$ws.Cells.Item(65, 1).Value = 'softball pants mens black'
$ws.Cells.Item(66, 1).Value = 'sport pants for men'
$ws.Cells.Item(67, 1).Value = 'basketball equipment pads'
$ws.Cells.Item(68, 1).Value = 'knee length leggings'
$ws.Cells.Item(69, 1).Value = 'volleyball knee pads girls youth'
$ws.Cells.Item(70, 1).Value = 'reduce swelling after surgery'
$ws.Cells.Item(71, 1).Value = 'soccer protection'
$ws.Cells.Item(72, 1).Value = 'knee support for basketball men'
$ws.Cells.Item(73, 1).Value = 'leg compression for men'
$ws.Cells.Item(74, 1).Value = 'knee support for basketball'
$ws.Cells.Item(75, 1).Value = 'soccer clothes'
$ws.Cells.Item(76, 1).Value = 'lacrosse equipment'
$ws.Cells.Item(77, 1).Value = 'capri athletic pants'
$ws.Cells.Item(78, 1).Value = 'basketball knee pads youth pair'
$ws.Cells.Item(79, 1).Value = 'mens leggings compression winter'
$ws.Cells.Item(80, 1).Value = 'ski compression pants men'
$ws.Cells.Item(81, 1).Value = 'leggings pants'
$ws.Cells.Item(82, 1).Value = 'youth volleyball knee pads'
$ws.Cells.Item(83, 1).Value = 'nike leggings basketball men'
$ws.Cells.Item(84, 1).Value = 'underarmour tights'
$ws.Cells.Item(85, 1).Value = 'yourh basketball knee pads'
$ws.Cells.Item(86, 1).Value = 'compression leggings knee pads'
$ws.Cells.Item(87, 1).Value = 'basketball tights with knee pads'
$ws.Cells.Item(88, 1).Value = 'leggings for men with knee pads'
$ws.Cells.Item(89, 1).Value = 'basketball knee pads men'
$ws.Cells.Item(90, 1).Value = 'mens compression capri pants'
$ws.Cells.Item(91, 1).Value = 'lacrosse leggings'
$ws.Cells.Item(92, 1).Value = 'mens compression leggings capri'
$ws.Cells.Item(93, 1).Value = 'basketball hip pads'
$ws.Cells.Item(94, 1).Value = 'male compression tights'
$ws.Cells.Item(95, 1).Value = 'compression tights for girls'
$ws.Cells.Item(96, 1).Value = 'basketball knee pads'
$ws.Cells.Item(97, 1).Value = 'mens compression pants capri'
$ws.Cells.Item(98, 1).Value = 'compression mens pants'
$ws.Cells.Item(99, 1).Value = 'basketball leggings for men'
$ws.Cells.Item(100, 1).Value = 'youth capri compression tights'
